$wb = $excel.ActiveWorkbook

# ======== Sheet: ALC ========
$ws = $wb.Worksheets.Item("ALC")
# -- row 2 (hunk 0) --
$ws.Range("H2").Value = 900
$ws.Range("I2").Value = 900
$ws.Range("K2").Value = 900
$ws.Range("M2").Value = -787
# -- row 28 (hunk 1) --
$ws.Range("H28").Value = 2069.1667
$ws.Range("I28").Value = 2816.25
$ws.Range("J28").Value = 575
$ws.Range("K28").Value = 2816.25
$ws.Range("L28").Value = 575
$ws.Range("M28").Value = -2331.25
$ws.Range("N28").Value = -1545
# -- row 40 (hunk 2) --
$ws.Range("H40").Value = 3626.7334
$ws.Range("I40").Value = 2340.1
$ws.Range("J40").Value = 6200
$ws.Range("K40").Value = 2340.1
$ws.Range("L40").Value = 6200
$ws.Range("M40").Value = -2165.1
$ws.Range("N40").Value = -6550
# -- row 43 (hunk 3) --
$ws.Range("H43").Value = 1000
$ws.Range("J43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("N43").Value = -1138
# -- row 48 (hunk 4) --
$ws.Range("H48").Value = 1000000
$ws.Range("J48").Value = 1000000
$ws.Range("L48").Value = 3000000
$ws.Range("N48").Value = -3000584
# -- row 56 (hunk 5) --
$ws.Range("H56").Value = 1000000
$ws.Range("J56").Value = 1000000
$ws.Range("L56").Value = 3000000
$ws.Range("N56").Value = -3001068
# -- row 62 (hunk 6) --
$ws.Range("H62").Value = 8329.5
$ws.Range("I62").Value = 8329.5
$ws.Range("K62").Value = 8329.5
$ws.Range("M62").Value = -7705.5
# -- row 65 (hunk 7) --
$ws.Range("H65").Value = 8329.5
$ws.Range("I65").Value = 8329.5
$ws.Range("K65").Value = 41647.5
$ws.Range("M65").Value = -38527.5
# -- row 111 (hunk 8) --
$ws.Range("H111").Value = 2717.1428
$ws.Range("J111").Value = 6016
$ws.Range("L111").Value = 18048
$ws.Range("N111").Value = -24182
# -- row 132 (hunk 9) --
$ws.Range("H132").Value = 1026.1765
$ws.Range("I132").Value = 1078.5333
$ws.Range("J132").Value = 633.5
$ws.Range("K132").Value = 3235.5999
$ws.Range("L132").Value = 1900.5
$ws.Range("M132").Value = -705.5999000000002
$ws.Range("N132").Value = -6960.5
# -- row 135 (hunk 10) --
$ws.Range("H135").Value = 280
$ws.Range("I135").Value = 280
$ws.Range("K135").Value = 2520
$ws.Range("M135").Value = 15
# -- row 137 (hunk 11) --
$ws.Range("H137").Value = 923.4167
$ws.Range("I137").Value = 923.4167
$ws.Range("K137").Value = 2770.2501
$ws.Range("M137").Value = -220.2501000000002
# -- row 138 (hunk 12) --
$ws.Range("H138").Value = 2286.875
$ws.Range("J138").Value = 2331.7827
$ws.Range("L138").Value = 6995.348100000001
$ws.Range("N138").Value = -17275.3481
# -- row 140 (hunk 13) --
$ws.Range("H140").Value = 150000
$ws.Range("I140").Value = 150000
$ws.Range("K140").Value = 150000
$ws.Range("M140").Value = -144820

# ======== Sheet: ARM ========
$ws = $wb.Worksheets.Item("ARM")
# -- row 5 (hunk 14) --
$ws.Range("H5").Value = 80.666664
$ws.Range("I5").Value = 50
$ws.Range("K5").Value = 50
$ws.Range("M5").Value = 62
# -- row 61 (hunk 15) --
$ws.Range("H61").Value = 1899.5
$ws.Range("I61").Value = 1899.5
$ws.Range("K61").Value = 1899.5
$ws.Range("M61").Value = -1687.5
# -- row 136 (hunk 16) --
$ws.Range("H136").Value = 1899.5
$ws.Range("I136").Value = 1899.5
$ws.Range("K136").Value = 5698.5
$ws.Range("M136").Value = -3148.5

# ======== Sheet: BSM ========
$ws = $wb.Worksheets.Item("BSM")
# -- row 4 (hunk 17) --
$ws.Range("H4").Value = 80.666664
$ws.Range("I4").Value = 50
$ws.Range("K4").Value = 50
$ws.Range("M4").Value = 65
# -- row 134 (hunk 18) --
$ws.Range("H134").Value = 2433.1667
$ws.Range("I134").Value = 2359.8
$ws.Range("K134").Value = 7079.400000000001
$ws.Range("M134").Value = -4544.400000000001

# ======== Sheet: CRP ========
$ws = $wb.Worksheets.Item("CRP")
# -- row 4 (hunk 19) --
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# -- row 16 (hunk 20) --
$ws.Range("H16").Value = 901.8333
$ws.Range("I16").Value = 959.6
$ws.Range("J16").Value = 613
$ws.Range("K16").Value = 959.6
$ws.Range("L16").Value = 613
$ws.Range("M16").Value = -672.6
$ws.Range("N16").Value = -1187
# -- row 58 (hunk 21) --
$ws.Range("H58").Value = 2488
$ws.Range("I58").Value = 2750.6667
$ws.Range("K58").Value = 2750.6667
$ws.Range("M58").Value = -2547.6667
# -- row 99 (hunk 22) --
$ws.Range("H99").Value = 4199.273
$ws.Range("I99").Value = 4419.2
$ws.Range("K99").Value = 4419.2
$ws.Range("M99").Value = -2921.2
# -- row 106 (hunk 23) --
$ws.Range("H106").Value = 60000
$ws.Range("J106").Value = 60000
$ws.Range("L106").Value = 60000
$ws.Range("N106").Value = -62524
# -- row 113 (hunk 24) --
$ws.Range("H113").Value = 901.8333
$ws.Range("I113").Value = 959.6
$ws.Range("J113").Value = 613
$ws.Range("K113").Value = 959.6
$ws.Range("L113").Value = 613
$ws.Range("M113").Value = 1210.4
$ws.Range("N113").Value = -4953
# -- row 126 (hunk 25) --
$ws.Range("H126").Value = 4199.273
$ws.Range("I126").Value = 4419.2
$ws.Range("K126").Value = 13257.6
$ws.Range("M126").Value = -10787.6
# -- row 132 (hunk 26) --
$ws.Range("H132").Value = 2599.3845
$ws.Range("I132").Value = 1846.1111
$ws.Range("K132").Value = 5538.3333
$ws.Range("M132").Value = -3008.3333
# -- row 136 (hunk 27) --
$ws.Range("H136").Value = 2488
$ws.Range("I136").Value = 2750.6667
$ws.Range("K136").Value = 8252.000100000001
$ws.Range("M136").Value = -5702.000100000001
# -- row 141 (hunk 28) --
$ws.Range("H141").Value = 50995.5
$ws.Range("J141").Value = 50995.5
$ws.Range("L141").Value = 50995.5
$ws.Range("N141").Value = -61355.5

# ======== Sheet: CUL ========
$ws = $wb.Worksheets.Item("CUL")
# -- row 129 (hunk 29) --
$ws.Range("H129").Value = 1707.875
$ws.Range("I129").Value = 1287.5
$ws.Range("K129").Value = 3862.5
$ws.Range("M129").Value = 1137.5
# -- row 130 (hunk 30) --
$ws.Range("H130").Value = 1000
$ws.Range("I130").Value = 1000
$ws.Range("K130").Value = 3000
$ws.Range("M130").Value = 2020

# ======== Sheet: GSM ========
$ws = $wb.Worksheets.Item("GSM")
# -- row 5 (hunk 31) --
$ws.Range("H5").Value = 26000
$ws.Range("I5").Value = 25000
$ws.Range("K5").Value = 25000
$ws.Range("M5").Value = -24888
# -- row 107 (hunk 32) --
$ws.Range("H107").Value = 723.5
$ws.Range("I107").Value = 1366.3334
$ws.Range("J107").Value = 80.666664
$ws.Range("K107").Value = 1366.3334
$ws.Range("L107").Value = 80.666664
$ws.Range("M107").Value = 553.6666
$ws.Range("N107").Value = -3920.666664
# -- row 132 (hunk 33) --
$ws.Range("H132").Value = 3903.6924
$ws.Range("I132").Value = 3716.6667
$ws.Range("J132").Value = 4324.5
$ws.Range("K132").Value = 11150.0001
$ws.Range("L132").Value = 12973.5
$ws.Range("M132").Value = -8620.000100000001
$ws.Range("N132").Value = -18033.5

# ======== Sheet: LTW ========
$ws = $wb.Worksheets.Item("LTW")
# -- row 2 (hunk 34) --
$ws.Range("H2").Value = 700000000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
# -- row 61 (hunk 35) --
$ws.Range("H61").Value = 750
$ws.Range("I61").Value = 750
$ws.Range("K61").Value = 750
$ws.Range("M61").Value = -548
# -- row 100 (hunk 36) --
$ws.Range("H100").Value = 2159.375
$ws.Range("I100").Value = 2159.375
$ws.Range("K100").Value = 2159.375
$ws.Range("M100").Value = -1618.375
# -- row 113 (hunk 37) --
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420

# ======== Sheet: WVR ========
$ws = $wb.Worksheets.Item("WVR")
# -- row 2 (hunk 38) --
$ws.Range("H2").Value = 148588.67
$ws.Range("I2").Value = 181485.83
$ws.Range("J2").Value = 17000
$ws.Range("K2").Value = 181485.83
$ws.Range("L2").Value = 17000
$ws.Range("M2").Value = -181373.83
$ws.Range("N2").Value = -17224
# -- row 107 (hunk 39) --
$ws.Range("H107").Value = 1172.4117
$ws.Range("I107").Value = 978.1429000000001
$ws.Range("K107").Value = 2934.4287
$ws.Range("M107").Value = -1014.4287
# -- row 113 (hunk 40) --
$ws.Range("H113").Value = 1184.3
$ws.Range("I113").Value = 1366.4286
$ws.Range("J113").Value = 759.3333
$ws.Range("K113").Value = 4099.2858
$ws.Range("L113").Value = 2277.9999
$ws.Range("M113").Value = -1929.2858
$ws.Range("N113").Value = -6617.9999
# -- row 132 (hunk 41) --
$ws.Range("H132").Value = 4300.1665
$ws.Range("I132").Value = 3700.5
$ws.Range("K132").Value = 11101.5
$ws.Range("M132").Value = -8571.5
# -- row 136 (hunk 42) --
$ws.Range("H136").Value = 4598.8887
$ws.Range("I136").Value = 4598.8887
$ws.Range("K136").Value = 13796.6661
$ws.Range("M136").Value = -11246.6661
